# Re-order the data rows (2..20) of Sheet1 according to the permutation
# described by the commit diff. Row 1 (header) and row 2 stay untouched;
# rows 3-20 are re-shuffled, each new row taking on the *entire* former
# content (values) of another row - i.e. a pure row permutation, not a
# per-cell edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping: new row number -> source (old) row number
$rowMap = @{
    2  = 2
    3  = 6
    4  = 8
    5  = 10
    6  = 3
    7  = 4
    8  = 12
    9  = 5
    10 = 11
    11 = 9
    12 = 14
    13 = 7
    14 = 15
    15 = 19
    16 = 13
    17 = 18
    18 = 16
    19 = 20
    20 = 17
}

$lastCol = 47   # column AU
$lastRow = 20

# Read the full original block (rows 2..20, columns A..AU) into memory first,
# since we need the untouched source data before overwriting anything.
$srcRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, $lastCol))
$srcValues = $srcRange.Value()

# Build the new block in the desired order.
# NOTE: Range.Value() returns a 1-based [1..19, 1..47] array, but an array
# created via New-Object is 0-based [0..18, 0..46]. Keep indices straight.
$newValues = New-Object 'object[,]' 19, $lastCol

for ($newRow = 2; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    $destIdx = $newRow - 2   # 0-based index into the 19-row block
    $srcIdx  = $oldRow - 1   # 1-based index into the 19-row block (source array)
    for ($c = 1; $c -le $lastCol; $c++) {
        $newValues[$destIdx, ($c - 1)] = $srcValues[$srcIdx, $c]
    }
}

# Write the re-ordered block back in one shot.
$destRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, $lastCol))
$destRange.Value = $newValues
